$d = $word.ActiveDocument

# The changelog currently ends with a bulleted ("List Paragraph" / numId 2)
# item reading "Put version on General page". We need to append:
#   Beta 2.5.1                              (plain paragraph, no list)
#     - Couple bugfixes                     (bulleted, same list as before)
#     - Bit more mobile compatability       (bulleted, same list as before)

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# 1) Append "Couple bugfixes" right after the last existing paragraph.
#    InsertParagraphAfter() duplicates the paragraph formatting of the
#    paragraph it is called on, so this new paragraph automatically gets
#    the same "List Paragraph" style / numId 2 bullet as "Put version on
#    General page" did - exactly what we want for this line.
$lastPara.Range.InsertParagraphAfter()
$idxCouple = $count + 1
$pCouple = $d.Paragraphs.Item($idxCouple)
$pCouple.Range.Text = "Couple bugfixes"

# 2) Append "Bit more mobile compatability" after that, again inheriting
#    the same bullet/list formatting.
$pCouple.Range.InsertParagraphAfter()
$idxMobile = $idxCouple + 1
$pMobile = $d.Paragraphs.Item($idxMobile)
$pMobile.Range.Text = "Bit more mobile compatability"

# 3) Insert the "Beta 2.5.1" heading paragraph *before* "Couple bugfixes".
#    This new (still empty) paragraph takes over index $idxCouple, pushing
#    the two bullet paragraphs down by one.
$pCouple.Range.InsertParagraphBefore()
$idxBeta = $idxCouple
$pBeta = $d.Paragraphs.Item($idxBeta)

# Strip the inherited bullet/list formatting and paragraph style so this
# becomes a plain paragraph, matching the other "Beta X.X" headings used
# throughout the rest of the changelog.
$pBeta.Range.ListFormat.RemoveNumbers()
$pBeta.Style = $d.Styles.Item("Normal")
$pBeta.Range.Text = "Beta 2.5.1"
